# Adds a new service-event row (row 16) to the "Card1" sheet and fills the
# previously-blank tracking cells in rows 2-15 with the literal text "nan"
# (matching how the source data was produced/exported).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# --- Fill previously blank cells (rows 2-15) with literal text "nan" ---
$ws.Range("D2:K2").Value2 = "nan"
$ws.Range("N2").Value2 = "nan"
$ws.Range("G3:K3").Value2 = "nan"
$ws.Range("M3:O3").Value2 = "nan"
$ws.Range("D4:K4").Value2 = "nan"
$ws.Range("D5:K5").Value2 = "nan"
$ws.Range("D6:K6").Value2 = "nan"
$ws.Range("M6").Value2 = "nan"
$ws.Range("D7:O7").Value2 = "nan"
$ws.Range("D8").Value2 = "nan"
$ws.Range("H8").Value2 = "nan"
$ws.Range("J8:K8").Value2 = "nan"
$ws.Range("M8:O8").Value2 = "nan"
$ws.Range("E9:G9").Value2 = "nan"
$ws.Range("I9").Value2 = "nan"
$ws.Range("K9").Value2 = "nan"
$ws.Range("M9:O9").Value2 = "nan"
$ws.Range("E10").Value2 = "nan"
$ws.Range("G10:J10").Value2 = "nan"
$ws.Range("M10:O10").Value2 = "nan"
$ws.Range("E11:F11").Value2 = "nan"
$ws.Range("H11:K11").Value2 = "nan"
$ws.Range("M11:O11").Value2 = "nan"
$ws.Range("D12:O12").Value2 = "nan"
$ws.Range("D13:O13").Value2 = "nan"
$ws.Range("D14:O14").Value2 = "nan"
$ws.Range("D15:O15").Value2 = "nan"

# --- Append the new event as row 16 ---
# Column A holds numeric-looking codes stored as text throughout the sheet,
# so force a text number format before writing "1" to keep it a string.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value2 = "1"

$ws.Range("L16").Value2 = "30/6/2025"
$ws.Range("M16").Value2 = "قطع سير كويلر مسنن 1270"
$ws.Range("N16").Value2 = "تم تغير سير 1270"
$ws.Range("O16").Value2 = "فني"
